$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = "admin"
$ws.Range("B3").Value = "wrong123"

# Row 4
$ws.Range("A4").Value = "user1"
$ws.Range("B4").Value = "pass1"

# Row 5
$ws.Range("A5").Value = "test"
$ws.Range("B5").Value = 12345
$ws.Range("B5").HorizontalAlignment = -4131

# Row 6
$ws.Range("A6").Value = "admin123"
$ws.Range("B6").Value = "Admin"

$ws.Range("D10").Select()
